$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the cell being edited (matches resulting selection in the file)
$ws.Range("A11").Select()

# Replace "Assert" with "Then" to match BDD syntax
$ws.Range("A11").Value = "Then"
